# Update underlying data that feeds the "ga-overall-comparison" chart.
# commit message: "change reporter generator; add CONF4 to PSO"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (MGenClust++ averages): ARI / DB Index / Silhouette Coefficient
$ws.Range("B4").Value = 0.6695
$ws.Range("D4").Value = 0.6873
$ws.Range("F4").Value = 0.5791

# Row 7 (MGenClust++ +/- 1 S.D.): ARI / DB Index / Silhouette Coefficient
$ws.Range("B7").Value = 0.094
$ws.Range("D7").Value = 0.1473
$ws.Range("F7").Value = 0.0525

# Move the active selection (matches the recorded sheet view state after edit)
$ws.Range("D8").Select()
